$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2094.2593
$ws.Range("I15").Value = 2094.2593
$ws.Range("K15").Value = 6282.777900000001
$ws.Range("M15").Value = -6113.777900000001
$ws.Range("H33").Value = 673.61536
$ws.Range("I33").Value = 394
$ws.Range("K33").Value = 394
$ws.Range("M33").Value = -165
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H43").Value = 3458
$ws.Range("I43").Value = 3458
$ws.Range("K43").Value = 3458
$ws.Range("M43").Value = -3389
$ws.Range("H74").Value = 5055
$ws.Range("I74").Value = 6536
$ws.Range("J74").Value = 3785.5715
$ws.Range("K74").Value = 6536
$ws.Range("L74").Value = 3785.5715
$ws.Range("M74").Value = -5600
$ws.Range("N74").Value = -5657.5715
$ws.Range("H77").Value = 5055
$ws.Range("I77").Value = 6536
$ws.Range("J77").Value = 3785.5715
$ws.Range("K77").Value = 32680
$ws.Range("L77").Value = 18927.8575
$ws.Range("M77").Value = -28000
$ws.Range("N77").Value = -28287.8575
$ws.Range("H98").Value = 3425
$ws.Range("I98").Value = 3100
$ws.Range("J98").Value = 3750
$ws.Range("K98").Value = 3100
$ws.Range("L98").Value = 3750
$ws.Range("M98").Value = -1602
$ws.Range("N98").Value = -6746
$ws.Range("H122").Value = 3425
$ws.Range("I122").Value = 3100
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 9300
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -6850
$ws.Range("N122").Value = -16150
$ws.Range("H137").Value = 66135.69500000001
$ws.Range("I137").Value = 138665
$ws.Range("J137").Value = 3967.7144
$ws.Range("K137").Value = 415995
$ws.Range("L137").Value = 11903.1432
$ws.Range("M137").Value = -413445
$ws.Range("N137").Value = -17003.1432
$ws.Range("H141").Value = 5537.1665
$ws.Range("I141").Value = 4526.7393
$ws.Range("K141").Value = 13580.2179
$ws.Range("M141").Value = -8400.2179
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1995.2903
$ws.Range("I74").Value = 2038.0454
$ws.Range("K74").Value = 2038.0454
$ws.Range("M74").Value = -1164.0454
$ws.Range("H77").Value = 1995.2903
$ws.Range("I77").Value = 2038.0454
$ws.Range("K77").Value = 10190.227
$ws.Range("M77").Value = -5822.226999999999
$ws.Range("H117").Value = 61974.75
$ws.Range("J117").Value = 61974.75
$ws.Range("L117").Value = 61974.75
$ws.Range("N117").Value = -71152.75
$ws.Range("H122").Value = 5160.5127
$ws.Range("I122").Value = 4628.4165
$ws.Range("K122").Value = 13885.2495
$ws.Range("M122").Value = -11435.2495
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 53435.4
$ws.Range("I20").Value = 127639.25
$ws.Range("K20").Value = 127639.25
$ws.Range("M20").Value = -127392.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 73659
$ws.Range("J87").Value = 73659
$ws.Range("L87").Value = 73659
$ws.Range("N87").Value = -76031
$ws.Range("H90").Value = 73659
$ws.Range("J90").Value = 73659
$ws.Range("L90").Value = 220977
$ws.Range("N90").Value = -232833
$ws.Range("H114").Value = 16250
$ws.Range("J114").Value = 16250
$ws.Range("L114").Value = 16250
$ws.Range("N114").Value = -24928
$ws.Range("H132").Value = 2823.6924
$ws.Range("I132").Value = 2186.476
$ws.Range("K132").Value = 6559.428
$ws.Range("M132").Value = -4029.428
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2232.6
$ws.Range("J21").Value = 2495.75
$ws.Range("L21").Value = 7487.25
$ws.Range("N21").Value = -7833.25
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8685
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -7908
$ws.Range("H87").Value = 999.5
$ws.Range("I87").Value = 999.5
$ws.Range("K87").Value = 2998.5
$ws.Range("M87").Value = -1750.5
$ws.Range("H90").Value = 999.5
$ws.Range("I90").Value = 999.5
$ws.Range("K90").Value = 8995.5
$ws.Range("M90").Value = -2755.5
$ws.Range("H103").Value = 425.44446
$ws.Range("I103").Value = 580.75
$ws.Range("J103").Value = 301.2
$ws.Range("K103").Value = 1742.25
$ws.Range("L103").Value = 903.5999999999999
$ws.Range("M103").Value = -863.25
$ws.Range("N103").Value = -2661.6
$ws.Range("H132").Value = 1102.5454
$ws.Range("I132").Value = 980.6667
$ws.Range("J132").Value = 1248.8
$ws.Range("K132").Value = 8826.0003
$ws.Range("L132").Value = 11239.2
$ws.Range("M132").Value = -6296.0003
$ws.Range("N132").Value = -16299.2
$ws.Range("H136").Value = 1484.1111
$ws.Range("I136").Value = 1484.1111
$ws.Range("K136").Value = 4452.3333
$ws.Range("M136").Value = 647.6666999999998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 37500
$ws.Range("I34").Value = 37500
$ws.Range("K34").Value = 37500
$ws.Range("M34").Value = -37232
$ws.Range("H39").Value = 105984.5
$ws.Range("J39").Value = 105984.5
$ws.Range("L39").Value = 105984.5
$ws.Range("N39").Value = -107048.5
$ws.Range("H70").Value = 299199.4
$ws.Range("J70").Value = 298999
$ws.Range("L70").Value = 298999
$ws.Range("N70").Value = -299539
$ws.Range("H73").Value = 299199.4
$ws.Range("J73").Value = 298999
$ws.Range("L73").Value = 298999
$ws.Range("N73").Value = -300871
$ws.Range("H76").Value = 37500
$ws.Range("I76").Value = 37500
$ws.Range("K76").Value = 37500
$ws.Range("M76").Value = -37185
$ws.Range("H79").Value = 37500
$ws.Range("I79").Value = 37500
$ws.Range("K79").Value = 37500
$ws.Range("M79").Value = -36408
$ws.Range("H102").Value = 1550.875
$ws.Range("I102").Value = 1456.1482
$ws.Range("K102").Value = 1456.1482
$ws.Range("M102").Value = 165.8517999999999
$ws.Range("H107").Value = 1986.2307
$ws.Range("I107").Value = 1787.25
$ws.Range("K107").Value = 1787.25
$ws.Range("M107").Value = 132.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3422.6155
$ws.Range("I7").Value = 3354.4
$ws.Range("K7").Value = 3354.4
$ws.Range("M7").Value = -3242.4
$ws.Range("H22").Value = 1718.8462
$ws.Range("I22").Value = 1586.8
$ws.Range("K22").Value = 1586.8
$ws.Range("M22").Value = -1291.8
$ws.Range("H27").Value = 1718.8462
$ws.Range("I27").Value = 1586.8
$ws.Range("K27").Value = 1586.8
$ws.Range("M27").Value = -1479.8
$ws.Range("H93").Value = 2560.125
$ws.Range("I93").Value = 2096.8333
$ws.Range("J93").Value = 3950
$ws.Range("K93").Value = 2096.8333
$ws.Range("L93").Value = 3950
$ws.Range("M93").Value = -848.8332999999998
$ws.Range("N93").Value = -6446
$ws.Range("H98").Value = 180000
$ws.Range("J98").Value = 180000
$ws.Range("L98").Value = 180000
$ws.Range("N98").Value = -185990
$ws.Range("H100").Value = 1738
$ws.Range("I100").Value = 1738
$ws.Range("K100").Value = 1738
$ws.Range("M100").Value = -1197
$ws.Range("H104").Value = 115000
$ws.Range("J104").Value = 115000
$ws.Range("L104").Value = 115000
$ws.Range("N104").Value = -121988
$ws.Range("H126").Value = 3422.6155
$ws.Range("I126").Value = 3354.4
$ws.Range("K126").Value = 10063.2
$ws.Range("M126").Value = -7593.200000000001
$ws.Range("H132").Value = 2001410.8
$ws.Range("J132").Value = 1933.3334
$ws.Range("L132").Value = 5800.0002
$ws.Range("N132").Value = -10860.0002
$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 114457
$ws.Range("J16").Value = 114457
$ws.Range("L16").Value = 114457
$ws.Range("N16").Value = -115041
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("H51").Value = 40198.8
$ws.Range("J51").Value = 59955
$ws.Range("L51").Value = 59955
$ws.Range("N51").Value = -60975
$ws.Range("H64").Value = 91743.75
$ws.Range("I64").Value = 66989
$ws.Range("K64").Value = 66989
$ws.Range("M64").Value = -66741
$ws.Range("H67").Value = 91743.75
$ws.Range("I67").Value = 66989
$ws.Range("K67").Value = 66989
$ws.Range("M67").Value = -66131
$ws.Range("H75").Value = 52999.668
$ws.Range("I75").Value = 18999
$ws.Range("K75").Value = 18999
$ws.Range("M75").Value = -18063
$ws.Range("H78").Value = 52999.668
$ws.Range("I78").Value = 18999
$ws.Range("K78").Value = 56997
$ws.Range("M78").Value = -52317
$ws.Range("H107").Value = 833.7778
$ws.Range("I107").Value = 769.6
$ws.Range("J107").Value = 914
$ws.Range("K107").Value = 2308.8
$ws.Range("L107").Value = 2742
$ws.Range("M107").Value = -388.8000000000002
$ws.Range("N107").Value = -6582
